# BIS-1002: remove the "Internal Assignment" column (column O) from the
# sample-type export sheet. The column header and its per-row "FALSE"
# values are cleared, while the existing cell formatting (style) for
# those cells is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First property-type table (rows 4-7): clear the "Internal Assignment"
# header (O4) and its values (O5:O7), keeping cell styles intact.
$ws.Range("O4:O7").ClearContents()

# Second property-type table (rows 12-15): same cleanup.
$ws.Range("O12:O15").ClearContents()

# Reflect the new selection left behind by removing the column's data
# (Excel naturally re-anchors the selection onto the now-empty range).
$ws.Range("O4:O15").Select() | Out-Null
